$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write new cell values in the exact order that reproduces the target
#     shared-string table layout. The engine rebuilds sharedStrings.xml at
#     save time: strings still referenced by some cell keep their original
#     relative slot, and brand-new strings are appended in first-write
#     order - so the order cells are touched below is deliberate. ---
$ws.Range("C2").Value = "url"
$ws.Range("C3").Value = "https://ui.cogmento.com/"
$ws.Range("A3").Value = "rahulscreencast9892@gmail.com"
$ws.Range("B3").Value = "Ra987456321@"
$ws.Range("D2").Value = "browser"
$ws.Range("D3").Value = "chrome"
$ws.Range("E2").Value = "firstName"
$ws.Range("F2").Value = "lastName"
$ws.Range("B2").Value = "password"
$ws.Range("G2").Value = "email"
$ws.Range("H2").Value = "description"
$ws.Range("E3").Value = "Sumeet"
$ws.Range("F3").Value = "Desai"
$ws.Range("G3").Value = "sumeet.desai@gmail.com"
$ws.Range("H3").Value = "Create a follow up activity"

# --- Hyperlinks: add the three new ones (call order controls the new
#     rId3/rId4/rId5 assignment and append order, matching the diff). ---
$ws.Hyperlinks.Add($ws.Range("C3"), "https://ui.cogmento.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Ra987456321@") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:sumeet.desai@gmail.com") | Out-Null

# --- Re-apply the workbook's two canonical cell styles so the new cells
#     line up with the existing "normal" (bordered) and "hyperlink"
#     (bordered + hyperlink font) styles instead of inventing new ones.
#     Must run after Hyperlinks.Add, which otherwise stamps its own
#     one-off style. PasteSpecial only takes a single-area destination
#     at a time in this host, so it is called once per cell. ---
foreach ($addr in @("B2","C2","D2","E2","F2","G2","H2","D3","E3","F3","H3")) {
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

foreach ($addr in @("B3","C3","G3")) {
    $ws.Range("A3").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# --- Column widths for the newly populated columns ---
$ws.Columns.Item(3).ColumnWidth = 23.7
$ws.Columns.Item(7).ColumnWidth = 23.6
$ws.Columns.Item(8).ColumnWidth = 23.6

# --- Selection, matching the edited workbook's last active cell ---
$ws.Range("G11").Select() | Out-Null
